$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 405 (shifts existing
# rows 405:437 down to 407:439), mirroring how the source data was
# re-exported with a newer day's prices inserted into the weekly series.
$ws.Rows("405:406").Insert()

# New row 405: Coliflor, Primera, Vega Monumental Concepción, 2023-06-29
$ws.Cells.Item(405, 1).Value = 11
$ws.Cells.Item(405, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(405, 3).Value = "Bíobío"
$ws.Cells.Item(405, 4).Value = "2023-06-29"
$ws.Cells.Item(405, 5).Value = 8
$ws.Cells.Item(405, 6).Value = 100112008
$ws.Cells.Item(405, 7).Value = "Coliflor"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value = 2000
$ws.Cells.Item(405, 11).Value = 700
$ws.Cells.Item(405, 12).Value = 800
$ws.Cells.Item(405, 13).Value = 750
$ws.Cells.Item(405, 14).Value = '$/unidad'
$ws.Cells.Item(405, 15).Value = "Región Metropolitana"
$ws.Cells.Item(405, 16).Value = 750
$ws.Cells.Item(405, 17).Value = 1
$ws.Cells.Item(405, 18).Value = "Hortaliza"

# New row 406: Coliflor, Segunda, Vega Monumental Concepción, 2023-06-29
$ws.Cells.Item(406, 1).Value = 11
$ws.Cells.Item(406, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(406, 3).Value = "Bíobío"
$ws.Cells.Item(406, 4).Value = "2023-06-29"
$ws.Cells.Item(406, 5).Value = 8
$ws.Cells.Item(406, 6).Value = 100112008
$ws.Cells.Item(406, 7).Value = "Coliflor"
$ws.Cells.Item(406, 8).Value = "Sin especificar"
$ws.Cells.Item(406, 9).Value = "Segunda"
$ws.Cells.Item(406, 10).Value = 1000
$ws.Cells.Item(406, 11).Value = 600
$ws.Cells.Item(406, 12).Value = 600
$ws.Cells.Item(406, 13).Value = 600
$ws.Cells.Item(406, 14).Value = '$/unidad'
$ws.Cells.Item(406, 15).Value = "Región Metropolitana"
$ws.Cells.Item(406, 16).Value = 600
$ws.Cells.Item(406, 17).Value = 1
$ws.Cells.Item(406, 18).Value = "Hortaliza"
